# Update all SSNs (column B numeric values, and the shared-string SSN on
# row 14) to start with a leading 9 instead of their original leading digit,
# since that implies they are fake test values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of cell -> new value (first digit replaced with 9)
$updates = @{
    "B2"  = "946412419"
    "B4"  = "909360961"
    "B5"  = "987777434"
    "B6"  = "973351423"
    "B7"  = "909175891"
    "B8"  = "909175891"
    "B11" = "965507"
    "B12" = "969005507"
    "B13" = "969860000"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = [double]$updates[$addr]
}

# B14 holds a shared string "44641241A" (not a pure number, has a trailing
# letter) - update it the same way, keeping it text.
$ws.Range("B14").Value = "94641241A"
